$wb = $excel.ActiveWorkbook

# Rename sheets (tab names)
$wb.Worksheets.Item(1).Name = "GNG_TO-1651168750210655"
$wb.Worksheets.Item(2).Name = "NB_TO-16511687525401354"
$wb.Worksheets.Item(3).Name = "RS_TO-16511687525410497"
$wb.Worksheets.Item(4).Name = "TOL_TO-16511687525857453"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16511687526593652"

# Sheet 1 (GNG)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16511687501830645.csv"
$ws1.Range("B3").Value = "GNG_stims-16511687501952503.csv"
$ws1.Range("B4").Value = "go_stims-16511687501972446.csv"
$ws1.Range("B5").Value = "GNG_stims-16511687502096553.csv"

# Sheet 2 (NB)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "ZB-match_0-16511687506021879.csv"
$ws2.Range("B3").Value = "OB-16511687515234997.csv"
$ws2.Range("B4").Value = "ZB-match_3-16511687503139899.csv"
$ws2.Range("B5").Value = "TB-16511687525263343.csv"
$ws2.Range("B6").Value = "OB-16511687515929728.csv"
$ws2.Range("B7").Value = "OB-16511687506836762.csv"
$ws2.Range("B8").Value = "TB-16511687519255512.csv"
$ws2.Range("B9").Value = "ZB-match_3-16511687502781775.csv"
$ws2.Range("B10").Value = "TB-16511687520461123.csv"

# Sheet 4 (TOL)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16511687525554018.csv"
$ws4.Range("B3").Value = "ZM_stims-1651168752543047.csv"
$ws4.Range("B4").Value = "MM_stims-16511687525703924.csv"
$ws4.Range("B5").Value = "ZM_stims-16511687525563948.csv"
$ws4.Range("B6").Value = "MM_stims-16511687525857453.csv"
$ws4.Range("B7").Value = "ZM_stims-16511687525713966.csv"

# Sheet 5 (vSAT)
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16511687526015859.csv"
$ws5.Range("B3").Value = "SAT_stims-16511687525906227.csv"
$ws5.Range("B4").Value = "vSAT_stims-1651168752644145.csv"
$ws5.Range("B5").Value = "vSAT_stims-165116875261691.csv"
